$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2,1).Value = "Última actualización: 16:52:37"

$ws.Cells.Item(3,1).Value = "Total filas: 35"

$ws.Cells.Item(5,1).Value = "Hora_Scrap"
$ws.Cells.Item(5,2).Value = "Hora_Llegada"
$ws.Cells.Item(5,3).Value = "Linea"
$ws.Cells.Item(5,4).Value = "Minutos"
$ws.Cells.Item(5,5).Value = "Parada"

$ws.Cells.Item(6,1).Value = "16:46:42"
$ws.Cells.Item(6,2).Value = "16:47"
$ws.Cells.Item(6,3).Value = "15_ABASTO"
$ws.Cells.Item(6,4).Value = 1
$ws.Cells.Item(6,5).Value = "LP1912"

$ws.Cells.Item(7,1).Value = "16:50:41"
$ws.Cells.Item(7,2).Value = "16:50"
$ws.Cells.Item(7,3).Value = "10_OLMOS"
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = "LP1912"

$ws.Cells.Item(8,1).Value = "16:50:41"
$ws.Cells.Item(8,2).Value = "16:51"
$ws.Cells.Item(8,3).Value = "15_ABASTO"
$ws.Cells.Item(8,4).Value = 1
$ws.Cells.Item(8,5).Value = "LP1912"

$ws.Cells.Item(9,1).Value = "16:52:37"
$ws.Cells.Item(9,2).Value = "16:53"
$ws.Cells.Item(9,3).Value = "15_ABASTO"
$ws.Cells.Item(9,4).Value = 1
$ws.Cells.Item(9,5).Value = "LP1912"

$ws.Cells.Item(10,1).Value = "16:46:42"
$ws.Cells.Item(10,2).Value = "16:53"
$ws.Cells.Item(10,3).Value = "10_OLMOS"
$ws.Cells.Item(10,4).Value = 7
$ws.Cells.Item(10,5).Value = "LP1912"

$ws.Cells.Item(11,1).Value = "16:46:42"
$ws.Cells.Item(11,2).Value = "16:56"
$ws.Cells.Item(11,3).Value = "215C_EL PATO"
$ws.Cells.Item(11,4).Value = 10
$ws.Cells.Item(11,5).Value = "LP1912"

$ws.Cells.Item(12,1).Value = "16:46:42"
$ws.Cells.Item(12,2).Value = "17:01"
$ws.Cells.Item(12,3).Value = "16_SANTA ANA"
$ws.Cells.Item(12,4).Value = 15
$ws.Cells.Item(12,5).Value = "LP1912"

$ws.Cells.Item(13,1).Value = "16:46:42"
$ws.Cells.Item(13,2).Value = "17:03"
$ws.Cells.Item(13,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(13,4).Value = 17
$ws.Cells.Item(13,5).Value = "LP1912"

$ws.Cells.Item(14,1).Value = "16:46:42"
$ws.Cells.Item(14,2).Value = "17:04"
$ws.Cells.Item(14,3).Value = "14_ABASTO"
$ws.Cells.Item(14,4).Value = 18
$ws.Cells.Item(14,5).Value = "LP1912"

$ws.Cells.Item(15,1).Value = "16:46:42"
$ws.Cells.Item(15,2).Value = "17:07"
$ws.Cells.Item(15,3).Value = "15_ABASTO"
$ws.Cells.Item(15,4).Value = 21
$ws.Cells.Item(15,5).Value = "LP1912"

$ws.Cells.Item(16,1).Value = "16:46:42"
$ws.Cells.Item(16,2).Value = "17:13"
$ws.Cells.Item(16,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(16,4).Value = 27
$ws.Cells.Item(16,5).Value = "LP1912"

$ws.Cells.Item(17,1).Value = "16:46:42"
$ws.Cells.Item(17,2).Value = "17:14"
$ws.Cells.Item(17,3).Value = "10_OLMOS"
$ws.Cells.Item(17,4).Value = 28
$ws.Cells.Item(17,5).Value = "LP1912"

$ws.Cells.Item(18,1).Value = "16:46:42"
$ws.Cells.Item(18,2).Value = "17:17"
$ws.Cells.Item(18,3).Value = "17_ROMERO"
$ws.Cells.Item(18,4).Value = 31
$ws.Cells.Item(18,5).Value = "LP1912"

$ws.Cells.Item(19,1).Value = "16:50:41"
$ws.Cells.Item(19,2).Value = "17:17"
$ws.Cells.Item(19,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(19,4).Value = 27
$ws.Cells.Item(19,5).Value = "LP1912"

$ws.Cells.Item(20,1).Value = "16:52:37"
$ws.Cells.Item(20,2).Value = "17:20"
$ws.Cells.Item(20,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(20,4).Value = 28
$ws.Cells.Item(20,5).Value = "LP1912"

$ws.Cells.Item(21,1).Value = "16:46:42"
$ws.Cells.Item(21,2).Value = "17:23"
$ws.Cells.Item(21,3).Value = "16_SANTA ANA"
$ws.Cells.Item(21,4).Value = 37
$ws.Cells.Item(21,5).Value = "LP1912"

$ws.Cells.Item(22,1).Value = "16:46:42"
$ws.Cells.Item(22,2).Value = "17:24"
$ws.Cells.Item(22,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(22,4).Value = 38
$ws.Cells.Item(22,5).Value = "LP1912"

$ws.Cells.Item(23,1).Value = "16:50:41"
$ws.Cells.Item(23,2).Value = "17:34"
$ws.Cells.Item(23,3).Value = "10_OLMOS"
$ws.Cells.Item(23,4).Value = 44
$ws.Cells.Item(23,5).Value = "LP1912"

$ws.Cells.Item(24,1).Value = "16:46:42"
$ws.Cells.Item(24,2).Value = "17:35"
$ws.Cells.Item(24,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(24,4).Value = 49
$ws.Cells.Item(24,5).Value = "LP1912"

$ws.Cells.Item(25,1).Value = "16:52:37"
$ws.Cells.Item(25,2).Value = "17:36"
$ws.Cells.Item(25,3).Value = "27_EL RETIRO"
$ws.Cells.Item(25,4).Value = 44
$ws.Cells.Item(25,5).Value = "LP1912"

$ws.Cells.Item(26,1).Value = "16:46:42"
$ws.Cells.Item(26,2).Value = "17:38"
$ws.Cells.Item(26,3).Value = "17X38_ROMERO"
$ws.Cells.Item(26,4).Value = 52
$ws.Cells.Item(26,5).Value = "LP1912"

$ws.Cells.Item(27,1).Value = "16:46:42"
$ws.Cells.Item(27,2).Value = "17:44"
$ws.Cells.Item(27,3).Value = "215B_EL PATO"
$ws.Cells.Item(27,4).Value = 58
$ws.Cells.Item(27,5).Value = "LP1912"

$ws.Cells.Item(28,1).Value = "16:50:41"
$ws.Cells.Item(28,2).Value = "17:47"
$ws.Cells.Item(28,3).Value = "16_SANTA ANA"
$ws.Cells.Item(28,4).Value = 57
$ws.Cells.Item(28,5).Value = "LP1912"

$ws.Cells.Item(29,1).Value = "16:46:42"
$ws.Cells.Item(29,2).Value = "17:48"
$ws.Cells.Item(29,3).Value = "27_EL RETIRO"
$ws.Cells.Item(29,4).Value = 62
$ws.Cells.Item(29,5).Value = "LP1912"

$ws.Cells.Item(30,1).Value = "16:50:41"
$ws.Cells.Item(30,2).Value = "17:49"
$ws.Cells.Item(30,3).Value = "27_EL RETIRO"
$ws.Cells.Item(30,4).Value = 59
$ws.Cells.Item(30,5).Value = "LP1912"

$ws.Cells.Item(31,1).Value = "16:46:42"
$ws.Cells.Item(31,2).Value = "17:50"
$ws.Cells.Item(31,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(31,4).Value = 64
$ws.Cells.Item(31,5).Value = "LP1912"

$ws.Cells.Item(32,1).Value = "16:52:37"
$ws.Cells.Item(32,2).Value = "17:51"
$ws.Cells.Item(32,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(32,4).Value = 59
$ws.Cells.Item(32,5).Value = "LP1912"

$ws.Cells.Item(33,1).Value = "16:46:42"
$ws.Cells.Item(33,2).Value = "18:02"
$ws.Cells.Item(33,3).Value = "17_ROMERO"
$ws.Cells.Item(33,4).Value = 76
$ws.Cells.Item(33,5).Value = "LP1912"

$ws.Cells.Item(34,1).Value = "16:52:37"
$ws.Cells.Item(34,2).Value = "18:03"
$ws.Cells.Item(34,3).Value = "17_ROMERO"
$ws.Cells.Item(34,4).Value = 71
$ws.Cells.Item(34,5).Value = "LP1912"

$ws.Cells.Item(35,1).Value = "16:46:42"
$ws.Cells.Item(35,2).Value = "18:04"
$ws.Cells.Item(35,3).Value = "14_ABASTO"
$ws.Cells.Item(35,4).Value = 78
$ws.Cells.Item(35,5).Value = "LP1912"

$ws.Cells.Item(36,1).Value = "16:52:37"
$ws.Cells.Item(36,2).Value = "18:14"
$ws.Cells.Item(36,3).Value = "10_OLMOS"
$ws.Cells.Item(36,4).Value = 82
$ws.Cells.Item(36,5).Value = "LP1912"

$ws.Cells.Item(37,1).Value = "16:46:42"
$ws.Cells.Item(37,2).Value = "18:24"
$ws.Cells.Item(37,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(37,4).Value = 98
$ws.Cells.Item(37,5).Value = "LP1912"

$ws.Cells.Item(38,1).Value = "16:46:42"
$ws.Cells.Item(38,2).Value = "18:34"
$ws.Cells.Item(38,3).Value = "14X44_ABASTO"
$ws.Cells.Item(38,4).Value = 108
$ws.Cells.Item(38,5).Value = "LP1912"

$ws.Cells.Item(39,1).Value = "16:46:42"
$ws.Cells.Item(39,2).Value = "18:38"
$ws.Cells.Item(39,3).Value = "17X38_ROMERO"
$ws.Cells.Item(39,4).Value = 112
$ws.Cells.Item(39,5).Value = "LP1912"

$ws.Cells.Item(40,1).Value = "16:46:42"
$ws.Cells.Item(40,2).Value = "18:41"
$ws.Cells.Item(40,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(40,4).Value = 115
$ws.Cells.Item(40,5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(2,1).Value = "Última actualización: 16:52:37"

$ws.Cells.Item(3,1).Value = "Total filas: 4"

$ws.Cells.Item(5,1).Value = "Hora_Scrap"
$ws.Cells.Item(5,2).Value = "Hora_Llegada"
$ws.Cells.Item(5,3).Value = "Linea"
$ws.Cells.Item(5,4).Value = "Minutos"
$ws.Cells.Item(5,5).Value = "Parada"

$ws.Cells.Item(6,1).Value = "16:46:42"
$ws.Cells.Item(6,2).Value = "16:56"
$ws.Cells.Item(6,3).Value = "215C_EL PATO"
$ws.Cells.Item(6,4).Value = 10
$ws.Cells.Item(6,5).Value = "LP1912"

$ws.Cells.Item(7,1).Value = "16:46:42"
$ws.Cells.Item(7,2).Value = "17:44"
$ws.Cells.Item(7,3).Value = "215B_EL PATO"
$ws.Cells.Item(7,4).Value = 58
$ws.Cells.Item(7,5).Value = "LP1912"

$ws.Cells.Item(8,1).Value = "16:46:42"
$ws.Cells.Item(8,2).Value = "17:50"
$ws.Cells.Item(8,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(8,4).Value = 64
$ws.Cells.Item(8,5).Value = "LP1912"

$ws.Cells.Item(9,1).Value = "16:52:37"
$ws.Cells.Item(9,2).Value = "17:51"
$ws.Cells.Item(9,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(9,4).Value = 59
$ws.Cells.Item(9,5).Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)

$ws.Cells.Item(2,1).Value = "Última actualización: 16:52:37"

$ws.Cells.Item(3,1).Value = "Total filas: 4"

$ws.Cells.Item(5,1).Value = "Hora_Scrap"
$ws.Cells.Item(5,2).Value = "Hora_Llegada"
$ws.Cells.Item(5,3).Value = "Linea"
$ws.Cells.Item(5,4).Value = "Minutos"
$ws.Cells.Item(5,5).Value = "Parada"

$ws.Cells.Item(6,1).Value = "16:52:37"
$ws.Cells.Item(6,2).Value = "16:57"
$ws.Cells.Item(6,3).Value = "215C_LA PLATA"
$ws.Cells.Item(6,4).Value = 5
$ws.Cells.Item(6,5).Value = "L6203"

$ws.Cells.Item(7,1).Value = "16:50:41"
$ws.Cells.Item(7,2).Value = "16:58"
$ws.Cells.Item(7,3).Value = "215C_LA PLATA"
$ws.Cells.Item(7,4).Value = 8
$ws.Cells.Item(7,5).Value = "L6203"

$ws.Cells.Item(8,1).Value = "16:46:42"
$ws.Cells.Item(8,2).Value = "16:59"
$ws.Cells.Item(8,3).Value = "215C_LA PLATA"
$ws.Cells.Item(8,4).Value = 13
$ws.Cells.Item(8,5).Value = "L6203"

$ws.Cells.Item(9,1).Value = "16:46:42"
$ws.Cells.Item(9,2).Value = "18:21"
$ws.Cells.Item(9,3).Value = "215C_LA PLATA"
$ws.Cells.Item(9,4).Value = 95
$ws.Cells.Item(9,5).Value = "L6203"

